$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 21:52"

# Update Murcia's row (row 31) statistics
$ws.Range("B31").Value = 1487
$ws.Range("C31").Value = 476
$ws.Range("D31").Value = 905
$ws.Range("E31").Value = 106
